$d = $word.ActiveDocument
$lq = [char]0x201C
$rq = [char]0x201D

# -----------------------------------------------------------------
# 1) "X db_queries.get_open_bugs()" paragraph: merge the split runs
#    back together (removes the stray gramStart/gramEnd proofErr
#    markers left over from "bugs(" / ")" being separate runs).
# -----------------------------------------------------------------
$d.Content.Find.Execute("open_bugs", $false, $false, $false, $false, $false, `
    $true, 1, $false, "open_bugs", 2) | Out-Null
$d.Content.Find.Execute("()", $false, $false, $false, $false, $false, `
    $true, 1, $false, "()", 2) | Out-Null

# -----------------------------------------------------------------
# 2) Insert "X get_bug_details_by_id()" as a new bullet right after
#    "It should be possible to view the detail of individual bugs..."
# -----------------------------------------------------------------
$detailText = "It should be possible to view the detail of individual bugs"
$pDetail = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith($detailText)) {
        $pDetail = $cand
        break
    }
}
$rDetail = $pDetail.Range
$rDetail.Collapse(0)
$rDetail.InsertParagraphAfter() | Out-Null
$newIdx = $pDetail.Range.Information(1) # unused, placeholder to keep numbering sane
$newPara1 = $pDetail.Next()
$newPara1.Range.ListFormat.ListLevelNumber = 2
$newPara1.Range.Text = "X get_bug_details_by_id()"

# -----------------------------------------------------------------
# 3) Insert "X add_table_entry()" as a new bullet right after
#    "It should be possible to create bugs"
# -----------------------------------------------------------------
$createText = "It should be possible to create bugs"
$pCreate = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith($createText)) {
        $pCreate = $cand
        break
    }
}
$rCreate = $pCreate.Range
$rCreate.Collapse(0)
$rCreate.InsertParagraphAfter() | Out-Null
$newPara2 = $pCreate.Next()
$newPara2.Range.ListFormat.ListLevelNumber = 2
$newPara2.Range.Text = "X add_table_entry()"

# -----------------------------------------------------------------
# 4) "Bug has a <status> field with open/closed options" -> "modify_bug_by_id"
# -----------------------------------------------------------------
$find1 = "Bug has a " + $lq + "status" + $rq + " field with open/closed options"
$d.Content.Find.Execute($find1, $false, $false, $false, $false, $false, `
    $true, 1, $false, "modify_bug_by_id", 2) | Out-Null

# -----------------------------------------------------------------
# 5) "Top level db has a list of <users>" -> "X add_table_entry()"
# -----------------------------------------------------------------
$find2 = "Top level db has a list of " + $lq + "users" + $rq
$d.Content.Find.Execute($find2, $false, $false, $false, $false, $false, `
    $true, 1, $false, "X add_table_entry()", 2) | Out-Null
